# Applies the "Improvement unit test project" change:
#  - Removes the tab-selected flag from the "Persons" sheet
#  - Adds new person/car data to the "Persons2" sheet (row 8 + a styled,
#    empty row 9 cell) and gives it an explicit page setup
#  - Appends a new, empty "Persons3" sheet which becomes the active /
#    selected sheet in the workbook

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Persons")
$ws2 = $wb.Worksheets.Item("Persons2")

# ---------------------------------------------------------------------
# 1. Populate Persons2 (sheet2) row 8 with a duplicate "person" record
#    plus a second, nested "car" record in columns L:N.
# ---------------------------------------------------------------------

# Columns B:H mirror the existing Persons!B5:H5 record - copying
# (instead of re-typing) keeps the original cell styles (text / date /
# number formats) intact.
$ws1.Range("B5:H5").Copy($ws2.Range("B8"))

# L8:N8 hold a new "Ford demo" car (Name, Targa, BuildYear). Values are
# first produced as TEXT() formula results on a scratch cell, then
# pasted-as-values so the workbook's string table picks them up as
# plain shared strings (no formula left behind, no extra number format
# created) - matching the assignment order Targa, BuildYear, Name.
$scratch = $ws2.Range("Z1")

$scratch.Formula = "QY478AZ"
$scratch.Copy()
$ws2.Range("M8").PasteSpecial(-4163) # xlPasteValues

$scratch.Formula = '=TEXT(1999,"0")'
$scratch.Copy()
$ws2.Range("N8").PasteSpecial(-4163) # xlPasteValues

$scratch.Formula = "Ford demo"
$scratch.Copy()
$ws2.Range("L8").PasteSpecial(-4163) # xlPasteValues

$scratch.Clear()

$ws2.Range("H8").Value = 2009

# Row 9 only carries a (styled, empty) placeholder cell at L9 - give it
# the same "applied font" formatting as H5/H8 above so an (empty) cell
# record is retained on save instead of being dropped as blank.
$ws2.Range("L9").Font.Name = $ws1.Range("H5").Font.Name
$ws2.Range("L9").Font.Size = $ws1.Range("H5").Font.Size

# Give Persons2 an explicit page setup (it had none before).
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Selection ends on L9; Select() also (temporarily) activates Persons2.
$null = $ws2.Range("L9").Select()

# ---------------------------------------------------------------------
# 2. Persons no longer is the tab-selected sheet (Persons3 will be).
# ---------------------------------------------------------------------
# (handled implicitly: any sheet Activate()/Select() call moves
# tabSelected away from Persons automatically)

# ---------------------------------------------------------------------
# 3. Append a new, empty "Persons3" sheet after Persons2 and make it
#    the active / selected sheet (activeTab=2).
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $last)
$ws3.Name = "Persons3"
$ws3.Activate()
